# Update the "取得日時" (acquired datetime) column for rows 2-6 on the
# active sheet ("ランサーズ") from 2025-12-20 12:34:22 to 2025-12-20 12:45:02.
# These are plain text values (no date number formatting), so we write
# the new value as a string to each cell to avoid Excel auto-converting
# it into a date/time value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-12-20 12:45:02"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
